$wb = $excel.ActiveWorkbook

# --- Glyphs sheet: add a new row (g32 / et reliqua) ---
$glyphs = $wb.Worksheets.Item("Glyphs")
$glyphs.Range("A33").Value = "g32"
$glyphs.Range("B33").Value = "et reliqua"
$glyphs.Range("B33").Select()

# --- ana sheet: add a new row (name) and make it the active sheet ---
$ana = $wb.Worksheets.Item("ana")
$ana.Activate()
$ana.Range("A19").Value = "name"
$ana.Range("A19").Select()
